# Auto-generated PowerShell Excel COM-interop edit script
# Applies literal numeric cell updates per the commit diff (static data refresh; no formulas).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 68533.53
$ws.Range("I137").Value = 1900.4
$ws.Range("J137").Value = 101850.1
$ws.Range("K137").Value = 5701.200000000001
$ws.Range("L137").Value = 305550.3
$ws.Range("M137").Value = -3151.200000000001
$ws.Range("N137").Value = -310650.3
$ws.Range("H138").Value = 11630268
$ws.Range("I138").Value = 34483956
$ws.Range("J138").Value = 2952.3684
$ws.Range("K138").Value = 103451868
$ws.Range("L138").Value = 8857.1052
$ws.Range("M138").Value = -103446728
$ws.Range("N138").Value = -19137.1052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14561.902
$ws.Range("I32").Value = 16049.925
$ws.Range("K32").Value = 16049.925
$ws.Range("M32").Value = -15762.925
$ws.Range("H61").Value = 1001606.5
$ws.Range("I61").Value = 1060464.9
$ws.Range("J61").Value = 1014
$ws.Range("K61").Value = 1060464.9
$ws.Range("L61").Value = 1014
$ws.Range("M61").Value = -1060252.9
$ws.Range("N61").Value = -1438
$ws.Range("H132").Value = 20791.518
$ws.Range("I132").Value = 2198.682
$ws.Range("J132").Value = 102600
$ws.Range("K132").Value = 6596.045999999999
$ws.Range("L132").Value = 307800
$ws.Range("M132").Value = -4066.045999999999
$ws.Range("N132").Value = -312860
$ws.Range("H136").Value = 1001606.5
$ws.Range("I136").Value = 1060464.9
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 3181394.7
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = -3178844.7
$ws.Range("N136").Value = -8142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1543.1428
$ws.Range("I86").Value = 1374.1305
$ws.Range("K86").Value = 1374.1305
$ws.Range("M86").Value = -251.1305
$ws.Range("H89").Value = 1543.1428
$ws.Range("I89").Value = 1374.1305
$ws.Range("K89").Value = 6870.6525
$ws.Range("M89").Value = -1254.6525
$ws.Range("H99").Value = 1777.2727
$ws.Range("I99").Value = 1962.5
$ws.Range("J99").Value = 1671.4286
$ws.Range("K99").Value = 1962.5
$ws.Range("L99").Value = 1671.4286
$ws.Range("M99").Value = -464.5
$ws.Range("N99").Value = -4667.4286
$ws.Range("H134").Value = 28294.975
$ws.Range("I134").Value = 30264.838
$ws.Range("K134").Value = 90794.514
$ws.Range("M134").Value = -88259.514

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3107.7837
$ws.Range("I31").Value = 1699.6
$ws.Range("K31").Value = 1699.6
$ws.Range("M31").Value = -1404.6
$ws.Range("H34").Value = 3107.7837
$ws.Range("I34").Value = 1699.6
$ws.Range("K34").Value = 1699.6
$ws.Range("M34").Value = -1497.6
$ws.Range("H58").Value = 32989.625
$ws.Range("I58").Value = 1735.5
$ws.Range("J58").Value = 126752
$ws.Range("K58").Value = 1735.5
$ws.Range("L58").Value = 126752
$ws.Range("M58").Value = -1532.5
$ws.Range("N58").Value = -127158
$ws.Range("H86").Value = 9269593
$ws.Range("I86").Value = 3438.6
$ws.Range("J86").Value = 20852286
$ws.Range("K86").Value = 3438.6
$ws.Range("L86").Value = 20852286
$ws.Range("M86").Value = -2315.6
$ws.Range("N86").Value = -20854532
$ws.Range("H89").Value = 9269593
$ws.Range("I89").Value = 3438.6
$ws.Range("J89").Value = 20852286
$ws.Range("K89").Value = 17193
$ws.Range("L89").Value = 104261430
$ws.Range("M89").Value = -11577
$ws.Range("N89").Value = -104272662
$ws.Range("H99").Value = 21742652
$ws.Range("I99").Value = 3178.5715
$ws.Range("J99").Value = 55559612
$ws.Range("K99").Value = 3178.5715
$ws.Range("L99").Value = 55559612
$ws.Range("M99").Value = -1680.5715
$ws.Range("N99").Value = -55562608
$ws.Range("H120").Value = 6242.857
$ws.Range("I120").Value = 4980
$ws.Range("J120").Value = 9400
$ws.Range("K120").Value = 4980
$ws.Range("L120").Value = 9400
$ws.Range("M120").Value = -1351
$ws.Range("N120").Value = -16658
$ws.Range("H121").Value = 5984.1177
$ws.Range("I121").Value = 5542
$ws.Range("J121").Value = 9300
$ws.Range("K121").Value = 5542
$ws.Range("L121").Value = 9300
$ws.Range("M121").Value = -4232
$ws.Range("N121").Value = -11920
$ws.Range("H126").Value = 21742652
$ws.Range("I126").Value = 3178.5715
$ws.Range("J126").Value = 55559612
$ws.Range("K126").Value = 9535.7145
$ws.Range("L126").Value = 166678836
$ws.Range("M126").Value = -7065.7145
$ws.Range("N126").Value = -166683776
$ws.Range("H132").Value = 2227.3235
$ws.Range("I132").Value = 1658.0303
$ws.Range("K132").Value = 4974.090899999999
$ws.Range("M132").Value = -2444.090899999999
$ws.Range("H134").Value = 796.6875
$ws.Range("I134").Value = 796.6875
$ws.Range("K134").Value = 2390.0625
$ws.Range("M134").Value = 144.9375
$ws.Range("H136").Value = 32989.625
$ws.Range("I136").Value = 1735.5
$ws.Range("J136").Value = 126752
$ws.Range("K136").Value = 5206.5
$ws.Range("L136").Value = 380256
$ws.Range("M136").Value = -2656.5
$ws.Range("N136").Value = -385356

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H56").Value = 3311
$ws.Range("I56").Value = 3311
$ws.Range("K56").Value = 3311
$ws.Range("M56").Value = -2781
$ws.Range("H62").Value = 6018.6924
$ws.Range("J62").Value = 8233.375
$ws.Range("L62").Value = 24700.125
$ws.Range("N62").Value = -26072.125
$ws.Range("H65").Value = 6018.6924
$ws.Range("J65").Value = 8233.375
$ws.Range("L65").Value = 74100.375
$ws.Range("N65").Value = -80964.375
$ws.Range("H113").Value = 706.13043
$ws.Range("J113").Value = 854
$ws.Range("L113").Value = 2562
$ws.Range("N113").Value = -6902
$ws.Range("H114").Value = 2287.6667
$ws.Range("I114").Value = 2225.6
$ws.Range("K114").Value = 6676.799999999999
$ws.Range("M114").Value = -3422.799999999999
$ws.Range("H127").Value = 909
$ws.Range("J127").Value = 909
$ws.Range("L127").Value = 2727
$ws.Range("N127").Value = -12647
$ws.Range("H129").Value = 418844.34
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H131").Value = 693.9394
$ws.Range("J131").Value = 715.3333
$ws.Range("L131").Value = 2145.9999
$ws.Range("N131").Value = -12225.9999
$ws.Range("H137").Value = 25643950
$ws.Range("I137").Value = 1197
$ws.Range("J137").Value = 55560496
$ws.Range("K137").Value = 3591
$ws.Range("L137").Value = 166681488
$ws.Range("M137").Value = 1509
$ws.Range("N137").Value = -166691688

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3627.4092
$ws.Range("I80").Value = 2945.182
$ws.Range("J80").Value = 4309.636
$ws.Range("K80").Value = 2945.182
$ws.Range("L80").Value = 4309.636
$ws.Range("M80").Value = -1947.182
$ws.Range("N80").Value = -6305.636
$ws.Range("H83").Value = 3627.4092
$ws.Range("I83").Value = 2945.182
$ws.Range("J83").Value = 4309.636
$ws.Range("K83").Value = 14725.91
$ws.Range("L83").Value = 21548.18
$ws.Range("M83").Value = -9733.91
$ws.Range("N83").Value = -31532.18
$ws.Range("H100").Value = 42000
$ws.Range("J100").Value = 42000
$ws.Range("L100").Value = 42000
$ws.Range("N100").Value = -44164
$ws.Range("H126").Value = 4610.9355
$ws.Range("I126").Value = 3801.25
$ws.Range("K126").Value = 11403.75
$ws.Range("M126").Value = -8933.75
$ws.Range("H132").Value = 79078.64999999999
$ws.Range("I132").Value = 56924.895
$ws.Range("J132").Value = 500000
$ws.Range("K132").Value = 170774.685
$ws.Range("L132").Value = 1500000
$ws.Range("M132").Value = -168244.685
$ws.Range("N132").Value = -1505060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6890
$ws.Range("J61").Value = 8400
$ws.Range("L61").Value = 8400
$ws.Range("N61").Value = -8804
$ws.Range("H68").Value = 2398.9443
$ws.Range("I68").Value = 2309.0908
$ws.Range("J68").Value = 2540.1428
$ws.Range("K68").Value = 2309.0908
$ws.Range("L68").Value = 2540.1428
$ws.Range("M68").Value = -1560.0908
$ws.Range("N68").Value = -4038.1428
$ws.Range("H71").Value = 2398.9443
$ws.Range("I71").Value = 2309.0908
$ws.Range("J71").Value = 2540.1428
$ws.Range("K71").Value = 11545.454
$ws.Range("L71").Value = 12700.714
$ws.Range("M71").Value = -7801.454
$ws.Range("N71").Value = -20188.714
$ws.Range("H113").Value = 6890
$ws.Range("J113").Value = 8400
$ws.Range("L113").Value = 8400
$ws.Range("N113").Value = -12740
$ws.Range("H132").Value = 503376.62
$ws.Range("I132").Value = 575125.7
$ws.Range("J132").Value = 1133.3334
$ws.Range("K132").Value = 1725377.1
$ws.Range("L132").Value = 3400.0002
$ws.Range("M132").Value = -1722847.1
$ws.Range("N132").Value = -8460.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58824496
$ws.Range("I81").Value = 1064.2667
$ws.Range("K81").Value = 2128.5334
$ws.Range("M81").Value = -1067.5334
$ws.Range("H84").Value = 58824496
$ws.Range("I84").Value = 1064.2667
$ws.Range("K84").Value = 10642.667
$ws.Range("M84").Value = -5338.666999999999
$ws.Range("H96").Value = 1412.25
$ws.Range("I96").Value = 1400
$ws.Range("J96").Value = 1424.5
$ws.Range("K96").Value = 1400
$ws.Range("L96").Value = 1424.5
$ws.Range("M96").Value = -27
$ws.Range("N96").Value = -4170.5
$ws.Range("H126").Value = 1083.6154
$ws.Range("I126").Value = 1030.826
$ws.Range("K126").Value = 3092.478
$ws.Range("M126").Value = -622.4780000000001
$ws.Range("H132").Value = 549.70734
$ws.Range("I132").Value = 533.325
$ws.Range("J132").Value = 1205
$ws.Range("K132").Value = 1599.975
$ws.Range("L132").Value = 3615
$ws.Range("M132").Value = 930.0249999999999
$ws.Range("N132").Value = -8675
$ws.Range("H136").Value = 23462972
$ws.Range("I136").Value = 29494208
$ws.Range("K136").Value = 88482624
$ws.Range("M136").Value = -88480074

